$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = @("26.133.41", "  +0.12%  ")
    3  = @("1.668.45", "  -0.31%  ")
    4  = @($null, "  -0.27%  ")
    5  = @("210.27", "  -0.93%  ")
    6  = @("0.5202", "  -1.13%  ")
    7  = @("1.004", "  -0.20%  ")
    8  = @("0.2612", "  -2.27%  ")
    9  = @("0.06330", "  +0.59%  ")
    10 = @("21.10", $null)
    11 = @("0.07552", "  -0.61%  ")
    12 = @("1.673.90", "  -0.24%  ")
    13 = @("4.417", "  -1.79%  ")
    14 = @("0.5428", "  -4.45%  ")
    15 = @("0.000008029", "  -1.20%  ")
    16 = @("66.38", "  +1.08%  ")
    17 = @("26.176.43", "  +0.21%  ")
    18 = @($null, "  -0.23%  ")
    19 = @("4.735", "  -1.99%  ")
    20 = @("187.21", "  -0.70%  ")
    21 = @($null, "  -3.02%  ")
    22 = @("6.247", "  +0.98%  ")
    24 = @("150.04", "  +0.73%  ")
    25 = @("0.1235", "  -1.11%  ")
    26 = @("7.478", "  -1.97%  ")
    27 = @("15.73", "  -0.58%  ")
    28 = @("0.06300", "  -0.94%  ")
    29 = @("1.369", "  +1.03%  ")
    30 = @("1.281", "  -0.96%  ")
    31 = @("3.506", "  -0.59%  ")
    32 = @("3.417", "  -3.37%  ")
    33 = @("1.646", "  -0.87%  ")
    34 = @($null, "  -0.83%  ")
    35 = @("0.6005", "  -0.16%  ")
    36 = @("2.398", "  -0.90%  ")
    37 = @("2.764", "  +1.94%  ")
    38 = @("1.112.41", "  +1.98%  ")
    39 = @("0.01613", "  -0.64%  ")
    40 = @("6.054", "  -1.22%  ")
    41 = @("0.8633", "  -0.85%  ")
    42 = @($null, "  -0.13%  ")
    43 = @("100.69", "  +0.81%  ")
    44 = @("1.821.88", "  -0.26%  ")
    45 = @("0.00000000108", "  -1.71%  ")
    46 = @("55.43", "  -2.68%  ")
    47 = @("0.9998", $null)
    48 = @("8.026", "  +0.97%  ")
    49 = @($null, "  -0.01%  ")
    50 = @("0.4237", "  -0.91%  ")
    51 = @("5.897", "  -0.73%  ")
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $dVal = $vals[0]
    $eVal = $vals[1]
    if ($null -ne $dVal) {
        $ws.Cells.Item($row, 4).Value = $dVal
    }
    if ($null -ne $eVal) {
        $ws.Cells.Item($row, 5).Value = $eVal
    }
}
